# Update automàtic: dades i banners [2026-02-19 20:20]
# Refreshes DATA_EXTRACCIO timestamps and the latest observation values
# pulled from meteo.cat for each station row in Dades_Meteo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-19 20:18:43'
$ws.Range('I2').Value = '3.1 mm'
$ws.Range('E3').Value = '2026-02-19 20:18:45'
$ws.Range('I3').Value = '4.6 mm'
$ws.Range('E4').Value = '2026-02-19 20:18:48'
$ws.Range('J4').Value = '1009.8 hPa'
$ws.Range('K4').Value = '10.9 MJ/m2'
$ws.Range('E5').Value = '2026-02-19 20:18:50'
$ws.Range('G5').Value = '146 cm'
$ws.Range('I5').Value = '7.5 mm'
$ws.Range('E6').Value = '2026-02-19 20:18:53'
$ws.Range('J6').Value = '1009.9 hPa'
$ws.Range('E7').Value = '2026-02-19 20:18:55'
$ws.Range('J7').Value = '1010.9 hPa'
$ws.Range('E8').Value = '2026-02-19 20:18:58'
$ws.Range('J8').Value = '1010.6 hPa'
$ws.Range('E9').Value = '2026-02-19 20:19:00'
$ws.Range('O9').Value = '10.5 °C'
$ws.Range('E10').Value = '2026-02-19 20:19:03'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '70%'
$ws.Range('N10').Value = '4.5 °C 19:58 TU'
$ws.Range('O10').Value = '10.5 °C'
$ws.Range('E11').Value = '2026-02-19 20:19:06'
$ws.Range('E12').Value = '2026-02-19 20:19:07'
$ws.Range('E13').Value = '2026-02-19 20:19:08'
$ws.Range('J13').Value = '1011.2 hPa'
$ws.Range('L13').Value = '47.9 km/h - 186º 19:41 TU'
$ws.Range('E14').Value = '2026-02-19 20:19:09'
$ws.Range('E15').Value = '2026-02-19 20:19:10'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '75%'
$ws.Range('O15').Value = '9.9 °C'
$ws.Range('E16').Value = '2026-02-19 20:19:11'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '75%'
$ws.Range('I16').Value = '9.3 mm'
$ws.Range('E17').Value = '2026-02-19 20:19:12'
$ws.Range('H17').NumberFormat = '@'
$ws.Range('H17').Value = '81%'
$ws.Range('E18').Value = '2026-02-19 20:19:13'
$ws.Range('J18').Value = '1010.1 hPa'
$ws.Range('O18').Value = '11.8 °C'
$ws.Range('E19').Value = '2026-02-19 20:19:15'
$ws.Range('E20').Value = '2026-02-19 20:19:16'
$ws.Range('L20').Value = '88.9 km/h - 321º 19:45 TU'
$ws.Range('E21').Value = '2026-02-19 20:19:17'
$ws.Range('J21').Value = '1011.2 hPa'
$ws.Range('E22').Value = '2026-02-19 20:19:20'
$ws.Range('E23').Value = '2026-02-19 20:19:22'
$ws.Range('I23').Value = '9.6 mm'
$ws.Range('E24').Value = '2026-02-19 20:19:24'
$ws.Range('J24').Value = '1014.7 hPa'
$ws.Range('E25').Value = '2026-02-19 20:19:27'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '67%'
$ws.Range('I25').Value = '5.9 mm'
$ws.Range('E26').Value = '2026-02-19 20:19:29'
$ws.Range('J26').Value = '1009.8 hPa'
$ws.Range('L26').Value = '63.4 km/h - 316º 19:47 TU'
$ws.Range('O26').Value = '3.1 °C'
$ws.Range('E27').Value = '2026-02-19 20:19:32'
$ws.Range('L27').Value = '55.4 km/h - 244º 19:31 TU'
$ws.Range('E28').Value = '2026-02-19 20:19:35'
$ws.Range('J28').Value = '1009.7 hPa'
$ws.Range('E29').Value = '2026-02-19 20:19:38'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '76%'
$ws.Range('O29').Value = '10.7 °C'
$ws.Range('E30').Value = '2026-02-19 20:19:40'
$ws.Range('H30').NumberFormat = '@'
$ws.Range('H30').Value = '78%'
$ws.Range('J30').Value = '1009.9 hPa'
$ws.Range('E31').Value = '2026-02-19 20:19:43'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '50%'
$ws.Range('J31').Value = '1009.4 hPa'
$ws.Range('O31').Value = '11.8 °C'
$ws.Range('E32').Value = '2026-02-19 20:19:45'
$ws.Range('E33').Value = '2026-02-19 20:19:48'
$ws.Range('J33').Value = '1010.7 hPa'
$ws.Range('L33').Value = '58.7 km/h - 254º 19:50 TU'
$ws.Range('E34').Value = '2026-02-19 20:19:50'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '62%'
$ws.Range('E35').Value = '2026-02-19 20:19:53'
$ws.Range('J35').Value = '1016.2 hPa'
$ws.Range('E36').Value = '2026-02-19 20:19:55'
$ws.Range('J36').Value = '1010.2 hPa'
$ws.Range('E37').Value = '2026-02-19 20:19:58'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '68%'
$ws.Range('J37').Value = '1011.2 hPa'
$ws.Range('O37').Value = '6.0 °C'
$ws.Range('E38').Value = '2026-02-19 20:20:01'
$ws.Range('E39').Value = '2026-02-19 20:20:04'
$ws.Range('I39').Value = '4.7 mm'
$ws.Range('E40').Value = '2026-02-19 20:20:06'
$ws.Range('J40').Value = '1012.3 hPa'
$ws.Range('E41').Value = '2026-02-19 20:20:09'
$ws.Range('J41').Value = '1012.8 hPa'
$ws.Range('E42').Value = '2026-02-19 20:20:11'
$ws.Range('O42').Value = '11.3 °C'
$ws.Range('E43').Value = '2026-02-19 20:20:14'
$ws.Range('O43').Value = '9.1 °C'
$ws.Range('E44').Value = '2026-02-19 20:20:17'
$ws.Range('I44').Value = '9.0 mm'
$ws.Range('L44').Value = '63.7 km/h - 120º 19:35 TU'
$ws.Range('E45').Value = '2026-02-19 20:20:19'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '85%'
$ws.Range('J45').Value = '1015.3 hPa'
$ws.Range('E46').Value = '2026-02-19 20:20:22'
$ws.Range('J46').Value = '1015.6 hPa'
